$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.156892776489258
$ws.Range("B1").Value = 2.545237064361572
$ws.Range("C1").Value = 6.786661148071289
$ws.Range("D1").Value = 2.064767837524414
$ws.Range("E1").Value = 1.219468355178833
